# Refresh the cryptocurrency price/volume snapshot (GitHub Actions cron update)
# and fix a few rows where the coin ranking order shifted (rows 40-49).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.258.90"
$ws.Range("E2").Value = "  -6.06%  "
$ws.Range("D3").Value = "3.491.25"
$ws.Range("E3").Value = "  -2.23%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'389.14"
$ws.Range("E5").Value = "  -6.85%  "
$ws.Range("D6").Value = "'120.03"
$ws.Range("E6").Value = "  -7.71%  "
$ws.Range("D7").Value = "3.485.03"
$ws.Range("E7").Value = "  -2.13%  "
$ws.Range("D8").Value = "'0.584"
$ws.Range("E8").Value = "  -9.94%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "'0.672"
$ws.Range("E10").Value = "  -12.15%  "
$ws.Range("D11").Value = "'0.150"
$ws.Range("E11").Value = "  -15.33%  "
$ws.Range("D12").Value = "'0.0000325"
$ws.Range("E12").Value = "  -3.79%  "
$ws.Range("D13").Value = "'38.55"
$ws.Range("E13").Value = "  -8.83%  "
$ws.Range("D14").Value = "4.033.33"
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("D15").Value = "'9.13"
$ws.Range("E15").Value = "  -8.52%  "
$ws.Range("D16").Value = "'0.135"
$ws.Range("E16").Value = "  -3.37%  "
$ws.Range("D17").Value = "3.475.27"
$ws.Range("E17").Value = "  -3.10%  "
$ws.Range("D18").Value = "'18.66"
$ws.Range("E18").Value = "  -8.54%  "
$ws.Range("D19").Value = "'12.56"
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("D20").Value = "63.197.20"
$ws.Range("E20").Value = "  -6.09%  "
$ws.Range("D21").Value = "'1.01"
$ws.Range("E21").Value = "  -11.41%  "
$ws.Range("D22").Value = "'392.63"
$ws.Range("E22").Value = "  -14.93%  "
$ws.Range("D23").Value = "'13.85"
$ws.Range("E23").Value = "  +3.35%  "
$ws.Range("D24").Value = "'80.73"
$ws.Range("E24").Value = "  -8.40%  "
$ws.Range("D25").Value = "'2.84"
$ws.Range("E25").Value = "  -8.70%  "
$ws.Range("D26").Value = "'33.19"
$ws.Range("E26").Value = "  -5.76%  "
$ws.Range("D27").Value = "'5.14"
$ws.Range("E27").Value = "  +5.80%  "
$ws.Range("D28").Value = "'2.97"
$ws.Range("E28").Value = "  -11.63%  "
$ws.Range("D29").Value = "'8.70"
$ws.Range("E29").Value = "  -14.71%  "
$ws.Range("D30").Value = "'11.77"
$ws.Range("E30").Value = "  -5.18%  "
$ws.Range("D31").Value = "'2.52"
$ws.Range("E31").Value = "  -9.90%  "
$ws.Range("D32").Value = "'0.110"
$ws.Range("E32").Value = "  -7.02%  "
$ws.Range("E33").Value = "  -8.78%  "
$ws.Range("E34").Value = "  -7.74%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("D36").Value = "'36.54"
$ws.Range("E36").Value = "  -12.15%  "
$ws.Range("D37").Value = "'53.54"
$ws.Range("E37").Value = "  -5.76%  "
$ws.Range("D38").Value = "'0.0434"
$ws.Range("E38").Value = "  -11.87%  "
$ws.Range("D39").Value = "'0.993"
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").Value = "'2.69"
$ws.Range("E40").Value = "  +15.43%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "'0.130"
$ws.Range("E41").Value = "  -10.68%  "
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0₃0624"
$ws.Range("E42").Value = "  -12.03%  "
$ws.Range("D43").Value = "'3.07"
$ws.Range("E43").Value = "  +13.13%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'141.22"
$ws.Range("E44").Value = "  -4.98%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'2.72"
$ws.Range("E45").Value = "  -10.36%  "
$ws.Range("B46").Value = "LidoDAOToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D46").Value = "'3.05"
$ws.Range("E46").Value = "  -6.33%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'24.70"
$ws.Range("E47").Value = "  +14.39%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'1.93"
$ws.Range("E48").Value = "  -2.04%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "'2.44"
$ws.Range("E49").Value = "  -10.50%  "
$ws.Range("D50").Value = "'3.99"
$ws.Range("E50").Value = "  -7.50%  "
$ws.Range("D51").Value = "'0.274"
$ws.Range("E51").Value = "  -11.72%  "
